$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CB")

$ws.Range("B2").Value = 1841000000.0
$ws.Range("B3").Value = 26689000000.0
$ws.Range("B4").Value = 28530000000.0
$ws.Range("B6").Value = 122965000000.0
$ws.Range("B8").Value = 5443000000.0
$ws.Range("B9").Value = 13878000000.0
$ws.Range("B10").Value = 163447000000.0
$ws.Range("B12").Value = 113751000000.0
$ws.Range("B13").Value = 1405000000.0
$ws.Range("B14").Value = 115156000000.0
$ws.Range("B15").Value = 14879000000.0
$ws.Range("B16").Value = 482000000.0
$ws.Range("B18").Value = 308000000.0
$ws.Range("B20").Value = 17745000000.0
$ws.Range("B22").Value = 9318000000.0
$ws.Range("B23").Value = 11064000000.0
$ws.Range("B24").Value = 41637000000.0
$ws.Range("B25").Value = 3901000000.0
$ws.Range("B27").Value = 59076000000.0
$ws.Range("B30").Value = 449677000.0
$ws.Range("G32").Value = 14628000000.0
$ws.Range("G33").Value = 16274000000.0
